$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.4913160583263677
$ws.Range("D3").Value = 0.4292746354602025
$ws.Range("D4").Value = 0.4168292003343551
$ws.Range("D5").Value = 0.4987461688492617
$ws.Range("D6").Value = 0.4359617349308071
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0.4531438655149995
$ws.Range("D9").Value = 0.1717284294603882
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0.3434568589207765
$ws.Range("D12").Value = 0.4520293489365654
$ws.Range("D13").Value = 0.4372620042723136
$ws.Range("D14").Value = 0.284944738552986
$ws.Range("D15").Value = 0.6878424816569146
$ws.Range("D16").Value = 0
$ws.Range("D17").Value = 0.09928485186217145
$ws.Range("D18").Value = 0
$ws.Range("D19").Value = 0.3240456951797158
$ws.Range("D20").Value = 0.3404848147116188
$ws.Range("D21").Value = 0.1021640196897929
$ws.Range("D22").Value = 0.9145537289867187
$ws.Range("D23").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0.2679483607318659
$ws.Range("D26").Value = 0.1563109501253831
$ws.Range("D27").Value = 0.1935543791213895
$ws.Range("D28").Value = 0.375499210550757
$ws.Range("D29").Value = 0.2356273799572769
$ws.Range("D30").Value = 0.3435497353023126
$ws.Range("D31").Value = 0.6271013281322559
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 0.3439212408284573
$ws.Range("D34").Value = 0.3881303984396768
$ws.Range("D35").Value = 0.230890684498932
$ws.Range("D36").Value = 0.3908238135042259
$ws.Range("D37").Value = 0.3089068449893193
$ws.Range("D38").Value = 0
$ws.Range("D39").Value = 0.07820191325345964
$ws.Range("D40").Value = 0.3238599424166435
$ws.Range("D41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0.394538868765673
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0.4280672425002322
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0.3163369555122132
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0.413114145072908
$ws.Range("D54").Value = 0.1227825763908238
$ws.Range("D55").Value = 0
$ws.Range("D56").Value = 0.2204885297668803
$ws.Range("D57").Value = 0.2243893377913997
$ws.Range("D58").Value = 0.239899693507941
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0.4667966936008173
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0.2128726664809139
$ws.Range("D64").Value = 0.4276957369740875
$ws.Range("D65").Value = 0.2890312993405777
$ws.Range("D66").Value = 1
$ws.Range("D67").Value = 0
$ws.Range("D68").Value = 0.1637410606482771
$ws.Range("D69").Value = 0.1525958948639361
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0.1861242685984954
$ws.Range("D72").Value = 0.7430110522894029
$ws.Range("D73").Value = 0.5112844803566453
$ws.Range("D74").Value = 0.08386737252716635
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 0
$ws.Range("D77").Value = 0.1206464196154918
$ws.Range("D78").Value = 0.8964428345871644
$ws.Range("D79").Value = 0.1292839230983561
$ws.Range("D80").Value = 0.2642333054704189
$ws.Range("D81").Value = 0.1674561159097241
$ws.Range("D82").Value = 0
$ws.Range("D83").Value = 0.4392124082845733
